# Use the timezone/locale-aware dateTool.format(...) helper instead of the
# joda-time DateTime.toString(...) calls for the "Period" and position
# "Time" report cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Period:" value cell (row 6, column B)
$ws.Range("B6").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", from, locale, timezone)+" - "+dateTool.format("YYYY-MM-dd HH:mm:ss", to, locale, timezone)}'

# Position "Time" column template cell (row 9, column B)
$ws.Range("B9").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", position.fixTime, locale, timezone)}'

# Restore the workbook's last active selection to B2
$ws.Range("B2").Select()
